$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: new time-record entry (30.3.2020 / IO Control / Create necessary files) ---
# Insert a new row 11 (pushes nothing, since it's past the last used row) so that it
# inherits the same per-column formatting as row 10 immediately above it.
$ws.Rows("11:11").Insert(-4121) # xlShiftDown

$ws.Range("A11").Value = "30.3.2020"
$ws.Range("B11").Value = 0.42708333333333331
$ws.Range("C11").Value = 0.43402777777777773
$ws.Range("D11").Formula = "=C11-B11"
$ws.Range("E11").Value = "IO Control"
$ws.Range("F11").Value = "Create necessary files"

# --- Row 12: new blank entry row, only the Date/From columns pre-formatted ---
$ws.Rows("12:12").Insert(-4121) # xlShiftDown
$ws.Range("C12:F12").Clear()

$ws.Range("A12").Select()
